$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1 (2)")

# --- Update existing input values (columns C:H) for rows 5, 6, 7 and 13 ---
# Row 5
$ws.Range("C5").Value = 542.78499999999997
$ws.Range("D5").Value = 20.5443
$ws.Range("E5").Value = 6483.69
$ws.Range("F5").Value = 10.114000000000001
$ws.Range("G5").Value = 248.68299999999999
$ws.Range("H5").Value = 11.1846

# Row 6
$ws.Range("C6").Value = 497.96100000000001
$ws.Range("D6").Value = 48.900799999999997
$ws.Range("E6").Value = 5113.43
$ws.Range("F6").Value = 18.147300000000001
$ws.Range("G6").Value = 276.863
$ws.Range("H6").Value = 26.597100000000001

# Row 7
$ws.Range("C7").Value = 413.85
$ws.Range("D7").Value = 37.5167
$ws.Range("E7").Value = 4581.96
$ws.Range("F7").Value = 15.858700000000001
$ws.Range("G7").Value = 175.90199999999999
$ws.Range("H7").Value = 19.2317

# Row 13
$ws.Range("C13").Value = 326.32799999999997
$ws.Range("D13").Value = 22.5303
$ws.Range("E13").Value = 6132.27
$ws.Range("F13").Value = 14.306699999999999
$ws.Range("G13").Value = 231.851
$ws.Range("H13").Value = 22.470099999999999

# --- Append new rows 38, 39, 40 with additional measurement results ---
$ws.Range("A38").Value = "alternate for 10"
$ws.Range("C38").Value = 110.712
$ws.Range("D38").Value = 10.232200000000001
$ws.Range("E38").Value = 7375.32
$ws.Range("F38").Value = 21.662800000000001
$ws.Range("G38").Value = 229.49799999999999
$ws.Range("H38").Value = 27.5611

$ws.Range("A39").Value = "alternate for -30"
$ws.Range("C39").Value = 403.267
$ws.Range("D39").Value = 37.002800000000001
$ws.Range("E39").Value = 6156.92
$ws.Range("F39").Value = 16.875
$ws.Range("G39").Value = 311.79199999999997
$ws.Range("H39").Value = 26.801300000000001

$ws.Range("A40").Value = "alternate for -40"
$ws.Range("C40").Value = 464.40100000000001
$ws.Range("D40").Value = 31.343800000000002
$ws.Range("E40").Value = 5609.18
$ws.Range("F40").Value = 18.006799999999998
$ws.Range("G40").Value = 299.214
$ws.Range("H40").Value = 22.308700000000002

# --- Recalculate and update the active selection to reflect where the edits ended ---
$excel.CalculateFull()
$ws.Activate()
$ws.Range("H41").Select()

$wb.Save()
